$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.826.56'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.619.62'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.13'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.98'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.256'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -1.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0881'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = '1.851.32'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').Value = '1.615.67'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.551'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.45'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').Value = '27.829.03'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '225.77'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').Value = '0.0₃0711'
$ws.Range('E20').Value = '  -1.58%  '
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.30'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.89'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.06'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.89'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.17'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0478'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.99%  '
$ws.Range('D33').Value = '1.413.92'
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.05'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.59'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.971'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.14%  '
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.841'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.68%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '64.94'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.34'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.84%  '
$ws.Range('E45').Value = '  -3.92%  '
$ws.Range('D46').Value = '1.759.93'
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.13'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.38'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('E50').Value = '  -2.49%  '
$ws.Range('E51').Value = '  -0.54%  '
